$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "newMessage" column (G) - header + values for the 3 data rows
$ws.Range("G1").Value = "newMessage"
$ws.Range("G2").Value = "new"
$ws.Range("G3").Value = "new"
$ws.Range("G4").Value = "new"

# Give column G an explicit width, matching the author's manual resize
$ws.Columns.Item(7).ColumnWidth = 12.8

# Move the selection onto the newly added column, as in the authored workbook
$excel.Goto($ws.Range("G1:G4"))
